$wb = $excel.ActiveWorkbook

$wsLocal = $wb.Worksheets.Item("Local_setup")
$wsSchool = $wb.Worksheets.Item("SchoolManagementProject")

# --- Local_setup sheet (sheet1): new command line for D4, next to the
#     existing postgresql row ---
[void]$wsLocal.Activate()
$wsLocal.Range("D4").Value = "pg_ctl -D ^`"C^:^\PostGreSQL^\data^`" -l logfile start"
[void]$wsLocal.Range("D4").Select()

# --- SchoolManagementProject sheet (sheet3): append "openfeign" dependency
#     under the Student microservice block ---
[void]$wsSchool.Activate()
$wsSchool.Range("B22").Value = "openfeign"
[void]$wsSchool.Range("B23").Select()

# SchoolManagementProject ends up the active/visible tab when the workbook
# is saved.
[void]$wsSchool.Activate()
